# Add a new "2023" data column (Q) to the table, mirroring the existing
# per-year columns (D..P), and clear the stray S4 selection that was left
# over in the saved sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build column Q by copying the formatting of column P (the previous
#     last year's column) into Q for every populated row, then overwrite
#     the values with the new 2023 figures. Using copy/PasteSpecial
#     (formats only) keeps the new cells on the same shared cell styles
#     the workbook already uses instead of inventing new ones. ---
$srcRows = @(2, 3, 4, 5, 6)
foreach ($r in $srcRows) {
    $ws.Range("P$r").Copy()
    $ws.Range("Q$r").PasteSpecial(-4122)
}

# New 2023 figures (row 2 stays blank, like the rest of that header row)
$ws.Range("Q3").Value = 2023
$ws.Range("Q4").Value = 279.01945525291825
$ws.Range("Q5").Value = 1792.7
$ws.Range("Q6").Value = 6425

# The two data rows grew slightly taller to fit the extra column's wrapped
# text.
$ws.Rows.Item(4).RowHeight = 27
$ws.Rows.Item(5).RowHeight = 27.75

# Clear the leftover "S4" selection outside the table by resetting the
# active cell back to A1.
$ws.Range("A1").Select()
